$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Excel's ColumnWidth property (character units) differs from the stored
# OOXML <col width> value by a constant offset of 5/6 (0.8333333333333334)
# for this workbook's default font. Subtract that offset so the saved
# width matches the target integer value exactly.
$offset = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 43 - $offset
$ws.Columns.Item(4).ColumnWidth = 78 - $offset
$ws.Columns.Item(6).ColumnWidth = 16 - $offset
$ws.Columns.Item(8).ColumnWidth = 29 - $offset

# --- Row 2 updates ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1327607"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327607"
$ws.Range("C2").Value = "Project Coordinator"
$ws.Range("D2").Value = "Nugegoda, Sri Lanka"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Brand Corridor (Pvt) Ltd"

# --- Row 3 updates ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1324549"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1324549"
$ws.Range("C3").Value = "Sales Account Manager"
$ws.Range("D3").Value = "Nasr City, Al Manteqah Al Oula, Nasr City, Cairo Governorate 4450113, Egypt"
$ws.Range("F3").Value = "18 applicants"
$ws.Range("H3").Value = "M911 Marketing Emer-Agency"

# --- New row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1321215"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1321215"
$ws.Range("C4").Value = "Business Management and Analytics Intern"
$ws.Range("D4").Value = "Manipal, Karnataka, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "M.A.H.E."
